$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild row 1 with the new player-record layout (name/pos columns added,
# everything shifted right, the old URL column dropped, and a trailing
# empty marker + numeric 0 appended). Force text formatting so values like
# dates and numeric-looking strings are stored verbatim as text, not
# auto-converted by Excel -- then strip the resulting format so no style
# index lingers on the cell.
$ws.Range("A1:K1").NumberFormat = "@"

$ws.Range("A1").Value = "Blacknall"
$ws.Range("B1").Value = "Saeed"
$ws.Range("C1").Value = "WR"
$ws.Range("D1").Value = "2018-11-18"
$ws.Range("E1").Value = "10"
$ws.Range("F1").Value = "22.246"
$ws.Range("G1").Value = "OAK"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "ARI"
$ws.Range("J1").Value = "W 23-21"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = 0

$ws.Range("A1:K1").ClearFormats()
